$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Turn the existing plain-text URLs in H15 and H17 into real hyperlinks.
#    (Formatting is copied from H2, which already carries the workbook's
#    "Hyperlink" cell style, so we don't spawn a brand-new style entry.)
# ---------------------------------------------------------------------------
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Fill in the three new BOM rows (29, 30, 31).
# ---------------------------------------------------------------------------

# Row 29 - Perf Board/Screw Terminals
$ws.Range("B29").Value = "Perf Board/Screw Terminals"
$ws.Range("D29").Value = "Amazon"
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 15.99
$ws.Range("H29").Value = "https://www.amazon.com/gp/product/B07FFDCF22/"

# Row 30 - Threaded Inserts (text first; its URL is entered after row 31's)
$ws.Range("B30").Value = "Threaded Inserts "
$ws.Range("D30").Value = "Basement"
$ws.Range("F30").Value = 0
$ws.Rows.Item(30).RowHeight = 28.8

# Row 31 - Machine Screws
$ws.Range("B31").Value = "Machine Screws"
$ws.Range("D31").Value = "Basement"
$ws.Range("F31").Value = 0
$ws.Range("H31").Value = "www.amazon.com/gp/product/B07HVRJW5J/"
$ws.Rows.Item(31).RowHeight = 28.8

# H30's URL is entered last so the shared-string table ends up in the same
# order the original authors built it in.
$ws.Range("H30").Value = "https://www.amazon.com/gp/product/B07L96KVP3/"

# ---------------------------------------------------------------------------
# 3. Apply the wrap-text "Hyperlink" look to the new link cells.
# ---------------------------------------------------------------------------
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null
$ws.Range("H31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Register the actual hyperlinks (order matters for relationship ids).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("H22"), "https://www.adafruit.com/product/2046") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H17"), "https://www.adafruit.com/product/4468") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H15"), "https://www.digikey.com/en/products/detail/adam-tech/IEC-A-4/9832319") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H31"), "www.amazon.com/gp/product/B07HVRJW5J/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H30"), "https://www.amazon.com/gp/product/B07L96KVP3/") | Out-Null

# ---------------------------------------------------------------------------
# 5. Sheet view tweaks: narrower column C and a different scroll/selection.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 40.3
$ws.Range("A11").Select()
$ws.Range("D21").Select()

$wb.Save()
